$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, pushing existing rows 4.. down by one.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly record.
$ws.Cells.Item(4, 1).Value = 9
$ws.Cells.Item(4, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(4, 3).Value = "Metropolitana"
$ws.Cells.Item(4, 4).Value = [DateTime]"2022-07-08"
$ws.Cells.Item(4, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(4, 5).Value = 13
$ws.Cells.Item(4, 6).Value = 100112035
$ws.Cells.Item(4, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 34
$ws.Cells.Item(4, 11).Value = 20000
$ws.Cells.Item(4, 12).Value = 22000
$ws.Cells.Item(4, 13).Value = 21000
$ws.Cells.Item(4, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(4, 15).Value = "Hijuelas"
$ws.Cells.Item(4, 16).Value = 1400
$ws.Cells.Item(4, 17).Value = 15
$ws.Cells.Item(4, 18).Value = "Hortaliza"
